# Generate Report for Handoff
# - Updates the "Ready for handoff" rows (source files handed off but not yet
#   handed back) after a new handoff/generate cycle:
#     * Priority changes from "low" to "ht" for those rows
#     * Latest Handoff Datetime (zh-cn) is refreshed
#     * Latest HO Xliff Generate Date (Overview) / Latest Handoff Datetime (de-de)
#       is refreshed to the newer timestamp

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date" for the
# "Ready for handoff" rows (4-7) moves from 00:28:13 to 00:28:33
$wsOverview.Range("G4:G7").Value = "2016-08-16 00:28:33"

# zh-cn sheet: column E = "Priority" goes from "low" to "ht" and
# column H = "Latest Handoff Datetime" is refreshed for rows 4-7
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2016-08-16 00:28:28"

# de-de sheet: column E = "Priority" goes from "low" to "ht" and
# column H = "Latest Handoff Datetime" is refreshed for rows 4-7
$wsDe.Range("E4:E7").Value = "ht"
$wsDe.Range("H4:H7").Value = "2016-08-16 00:28:33"
